# Updated symbol list on Sun Dec 11 22:32:27 UTC 2022 with GitHub Actions
#
# Price refresh for the crypto table on Sheet1. Column D holds prices as
# text (e.g. "21.20"), so numeric-looking values must be written while the
# cell's NumberFormat is forced to Text ("@") -- otherwise Excel silently
# re-types the cell as a Number and drops significant trailing zeros
# (21.20 -> 21.2). After writing, the cell's format is restored from an
# untouched template cell (D13) so we don't leave a stray "Text" style
# behind on cells that were never styled before.
#
# Rows 41-43 additionally got re-ranked (BKEXToken/CEJI/KickToken rotated
# to new rows with refreshed prices), so those rows' Coin/Link/Volume
# columns are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D price updates (stored as text) ---------------------------
$priceUpdates = @{
    "D2"  = "286.48"
    "D3"  = "21.20"
    "D4"  = "6.452"
    "D5"  = "0.06375"
    "D6"  = "3.604"
    "D7"  = "1.553"
    "D8"  = "6.575"
    "D9"  = "0.8245"
    "D10" = "0.01420"
    "D11" = "0.1683"
    "D12" = "0.08769"
    "D14" = "0.03205"
    "D15" = "0.09196"
    "D16" = "3.702"
    "D17" = "0.001639"
    "D18" = "0.04747"
    "D19" = "0.006206"
    "D20" = "0.006291"
    "D24" = "2.322"
    "D25" = "0.3357"
    "D40" = "0.04793"
    "D41" = "0.1119"
    "D42" = "0.003454"
    "D43" = "0.007139"
    "D44" = "0.01171"
    "D45" = "0.00006978"
    "D48" = "0.006162"
    "D50" = "0.01241"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceUpdates[$addr]
}

# Restore the original (default) cell style on each updated price cell by
# copying formats from an untouched cell in the same column, so only the
# value changes -- not the style index.
$ws.Range("D13").Copy()
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Rows 41-43: Coin / Link / Volume(1h) re-rank -----------------------
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICK"
